$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = "Caitlin Boyd, poor"
    "B2" = "Lexi Green, poor"
    "C2" = "Stanley Hirst, poor"
    "D2" = "Ava Lee, poor"
    "E2" = "James Shilton, poor"
    "F2" = "James Calderon, poor"
    "G2" = "William Hunt, good"
    "H2" = "Violet Hudson, poor"
    "I2" = "Niko Morris, poor"
    "J2" = "Ruby Haigh, good"

    "A3" = ""
    "B3" = "Brooke Layton, good"
    "C3" = "Aarron Kelly, good"
    "D3" = "Katrina Petersone, good"
    "E3" = "Madison Taylor, good"
    "F3" = "Nancy Enyoazu, good"
    "G3" = "Benjamin Finn, good"
    "H3" = "Esther Sido, excellent"
    "I3" = "Benedict Hobday, good"
    "J3" = "Samuel Dixon, excellent"

    "A4" = ""
    "B4" = "Thomas Barrett, excellent"
    "C4" = "Benjamin Hillary, excellent"
    "D4" = "Spencer Rowe, excellent"
    "E4" = "James Eilbeck, excellent"
    "F4" = "Matthew Homan, excellent"
    "G4" = "Alex Sentance, excellent"
    "H4" = "Sophie Rayner, excellent"
    "I4" = ""
    "J4" = ""
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
